# Update database and change read_price algorithm
# Quarterly income statement: drop the oldest reporting period (column D,
# "6 ماهه منتهی به 1399/06") and append the newest period
# ("12 ماهه منتهی به 1401/12") as the new rightmost column (M), shifting
# everything else one column to the left.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Drop the oldest period: deleting column D shifts D:M left to C:L,
#    carrying values, shared-string text and per-cell styles with it, and
#    also shifts the <cols> width definitions left by one column.
$ws.Columns.Item(4).Delete()

# 2) Clone column L's formatting into the freshly vacated column M (row by
#    row the style of every data column in this sheet is identical, so
#    copying L's formats is exactly what the new column needs).
$ws.Range($ws.Cells.Item(1, 12), $ws.Cells.Item(28, 12)).Copy()
$ws.Range($ws.Cells.Item(1, 13), $ws.Cells.Item(28, 13)).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Column M is a "period end" column (like the old F/J columns), so it gets
# the wider 29-width treatment instead of the standard 28.
$ws.Columns.Item(13).ColumnWidth = 28.17

# 3) Fix up the header row text for the newly introduced period/date.
$ws.Cells.Item(8, 13).Value = "12 ماهه منتهی به 1401/12"
$ws.Cells.Item(9, 13).Value = "1402-02-25 (2)"

# The publish-date that used to read "1401-10-28 (8)" (column J) is now in
# column I after the shift, and the re-filed report carries a new date.
$ws.Cells.Item(9, 9).Value = "1402-02-25 (10)"

# 4) Fill in the new rightmost column with the latest quarter's figures.
$ws.Cells.Item(11, 13).Value = 56612
$ws.Cells.Item(12, 13).Value = -26760
$ws.Cells.Item(13, 13).Value = 29852
$ws.Cells.Item(14, 13).Value = -7422
$ws.Cells.Item(15, 13).Value = "-"
$ws.Cells.Item(16, 13).Value = 4920
$ws.Cells.Item(17, 13).Value = 27349
$ws.Cells.Item(18, 13).Value = -52
$ws.Cells.Item(19, 13).Value = 2391
$ws.Cells.Item(20, 13).Value = 29688
$ws.Cells.Item(21, 13).Value = -1165
$ws.Cells.Item(22, 13).Value = 28523
$ws.Cells.Item(23, 13).Value = "-"
$ws.Cells.Item(24, 13).Value = 28523
$ws.Cells.Item(25, 13).Value = 0
$ws.Cells.Item(26, 13).Value = 3967
$ws.Cells.Item(27, 13).Value = 0
